# Auto commit 30-05-2025 13:26:06.96
# Applies the "Cartage@ 1.00 %" label update + the re-derivation of the
# "Cartage" / "Contractors Profit & Overhead" rows so that they are based
# on the MATERIAL-only subtotal (F-column of the material section) instead
# of the previous running UNIT TOTAL, for every one of the 8 repeated
# pricing blocks in the sheet, plus the knock-on recalculation of the
# dependent TOTAL / UNIT TOTAL / "Say" / cross-check cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shared-string text tweak: "Cartage@ 1.00 %" -> "Cartage@ 1.00 % (for material)"
#    (every cell that currently holds the old label gets the new one)
# ---------------------------------------------------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Cartage@ 1.00 %") {
        $cell.Value2 = "Cartage@ 1.00 % (for material)"
    }
}

# ---------------------------------------------------------------------
# 2) Per-block formula surgery.
#    Each pricing block starts at row $b (8, 39, 70, 101, 132, 163, 194, 225)
#    and re-uses the same relative row layout:
#      b+0  .. material line            F(b)   = D(b)*E(b)
#      b+4  .. "Other charge" (M+L)     F(b+4) = (F(b)+F(b+2))*0.5%
#      b+5  .. TOTAL                    F(b+5) = F(b)+F(b+2)+F(b+4)
#      b+6  .. UNIT TOTAL               F(b+6) = F(b+5)/D(b+5)
#      b+7  .. Cartage (material)       E(b+7) [NEW] = F(b)*100.5%/D(b+5)
#                                        F(b+7)        = E(b+7)*1%      (was F(b+6)*1%)
#      b+8  .. TOTAL                    F(b+8) = F(b+6)+F(b+7)
#      b+9  .. Contractors P&O          E(b+9) [NEW] = F(b+8)
#                                        F(b+9)        = E(b+9)*15%    (was F(b+8)*15%)
#      b+10 .. UNIT TOTAL (final)       F(b+10) = round(F(b+8)+F(b+9),2)
#      b+12 .. labour line              F(b+12) = D(b+12)*E(b+12)
#      b+16 .. "Other charge" (M+L)     F(b+16) = (F(b+12)+F(b+14))*0.5%
#                                        H(b+16) [NEW] = F(b+16)+F(b+4)
#      b+17 .. TOTAL                    F(b+17) = F(b+12)+F(b+14)+F(b+16)
#      b+18 .. UNIT TOTAL               F(b+18) = F(b+17)/D(b+17)
#      b+19 .. Cartage (labour)         E(b+19) [NEW] = F(b+12)*100.5%/D(b+17)
#                                        F(b+19)        = E(b+19)*1%   (was F(b+18)*1%)
#      b+20 .. TOTAL                    F(b+20) = F(b+18)+F(b+19)
#      b+21 .. Contractors P&O          E(b+21) [NEW] = F(b+20)
#                                        F(b+21)        = E(b+21)*15% (was F(b+20)*15%)
#      b+22 .. UNIT TOTAL (final)       F(b+22) = round(F(b+20)+F(b+21),2)
#                                        I(b+22) [NEW] = F(b+22)+F(b+10)
#      b+23 .. "Say" row                H(b+23) and I(b+23) formulas change shape
#    For the last three blocks (163, 194, 225) the author additionally
#    cross-checks the material-cartage running total in column H at
#    b+19 and b+21.
# ---------------------------------------------------------------------

$blocks = @(8, 39, 70, 101, 132, 163, 194, 225)
$extraHBlocks = @(163, 194, 225)

foreach ($b in $blocks) {

    $r7  = $b + 7
    $r8  = $b + 8
    $r9  = $b + 9
    $r10 = $b + 10
    $r16 = $b + 16
    $r19 = $b + 19
    $r20 = $b + 20
    $r21 = $b + 21
    $r22 = $b + 22
    $r23 = $b + 23

    $matRow  = $b        # material line
    $unitRow = $b + 5    # row holding the qty divisor (D) for the material UNIT TOTAL
    $labRow  = $b + 12   # labour line
    $lUnitRow = $b + 17  # row holding the qty divisor (D) for the labour UNIT TOTAL
    $hRow    = $b + 4    # material "other charge" row (for H(b+16) cross add)
    $fin1Row = $b + 10   # material final UNIT TOTAL row (for I(b+22) add)

    # -- Cartage (material), row b+7 --------------------------------------------------
    $ws.Range("E$r7").Formula = "=F$matRow*100.5%/D$unitRow"
    $ws.Range("E$r7").NumberFormat = "0.00000"
    $ws.Range("F$r7").Formula = "=E$r7*1%"

    # -- Contractors P&O (material), row b+9 -------------------------------------------
    $ws.Range("E$r9").Formula = "=F$r8"
    $ws.Range("F$r9").Formula = "=E$r9*15%"

    # -- "Other charge" (labour), row b+16: add cross-check in column H ---------------
    $ws.Range("H$r16").Formula = "=F$r16+F$hRow"

    # -- Cartage (labour), row b+19 ----------------------------------------------------
    $ws.Range("E$r19").Formula = "=F$labRow*100.5%/D$lUnitRow"
    $ws.Range("E$r19").NumberFormat = "0.00000"
    $ws.Range("F$r19").Formula = "=E$r19*1%"

    # -- Contractors P&O (labour), row b+21 --------------------------------------------
    $ws.Range("E$r21").Formula = "=F$r20"
    $ws.Range("F$r21").Formula = "=E$r21*15%"

    # -- UNIT TOTAL final (labour), row b+22: add cross total in column I -------------
    $ws.Range("I$r22").Formula = "=F$r22+F$fin1Row"

    # -- "Say" row, b+23: H/I cross-check formulas change shape ------------------------
    $h1 = $b
    $h2 = $b + 2
    $i1 = $b + 12
    $i2 = $b + 14
    $ws.Range("H$r23").Formula = "=((F$h1*1.005*1.01)+(F$h2*1.005))*1.15"
    $ws.Range("I$r23").Formula = "=((F$i1*1.005*1.01)+(F$i2*1.005))*1.15"

    # -- Extra material-cartage cross checks in column H for blocks 6,7,8 -------------
    if ($extraHBlocks -contains $b) {
        $ws.Range("H$r19").Formula = "=E$r19+E$r7"
        $ws.Range("H$r21").Formula = "=E$r21+E$r9"
    }
}
